$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5

# Row 4 updates
$ws.Range("G4").Value = 1.48
$ws.Range("H4").Value = 4.2
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.95
$ws.Range("AE4").Value = 17
$ws.Range("AU4").Value = 8.5
$ws.Range("BA4").Value = 126
